$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily NAV rows for 2024-09-02 .. 2024-09-27 appended after the existing
# history (rows 3-650), continuing the basket/NAV series.
$newRows = @(
    ,@("2024-09-02", 1660.599975585938, 1280.449951171875, 1769.650024414062, 2220.5, 464.8500061035156, 8325.749969482422, 0, 240.8210630544845)
    ,@("2024-09-03", 1683.75, 1248.449951171875, 1810.949951171875, 2205.5, 469.6000061035156, 8357.449920654297, 0.003807458942205738, 241.7379793644828)
    ,@("2024-09-04", 1699.5, 1288.900024414062, 1780.25, 2217.949951171875, 466.5499877929688, 8386.249938964844, 0.003446029420932761, 242.5710155535297)
    ,@("2024-09-05", 1694.699951171875, 1272.300048828125, 1783.150024414062, 2236.14990234375, 467.75, 8389.549926757812, 0.0003934998142180441, 242.6664672030847)
    ,@("2024-09-06", 1724.449951171875, 1234.300048828125, 1778.650024414062, 2224.39990234375, 463.3999938964844, 8351.999908447266, -0.004475808432915339, 241.5803385827913)
    ,@("2024-09-09", 1789.300048828125, 1234.300048828125, 1817.949951171875, 2138.60009765625, 460.7999877929688, 8362.550109863281, 0.001263194627833399, 241.8855015686793)
    ,@("2024-09-10", 1799.949951171875, 1249.650024414062, 1830.099975585938, 2140.14990234375, 466.9500122070312, 8420.699890136719, 0.006953594239734628, 243.5674751990626)
    ,@("2024-09-11", 1788.25, 1255, 1828.699951171875, 2119.85009765625, 459.2999877929688, 8369.700012207031, -0.006056489198650145, 242.092311416377)
    ,@("2024-09-12", 1816.650024414062, 1257.449951171875, 1852.949951171875, 2111, 459, 8415.049926757812, 0.005418344084571652, 243.4040508598601)
    ,@("2024-09-13", 1888, 1264.349975585938, 1921.550048828125, 2101.35009765625, 451.75, 8530.500122070312, 0.01371949023681921, 246.7434303592342)
    ,@("2024-09-16", 1934.900024414062, 1259.75, 1916, 2047.199951171875, 447.8999938964844, 8501.549957275391, -0.003393724210849177, 245.9060512057561)
    ,@("2024-09-17", 1891.199951171875, 1257.550048828125, 1904.050048828125, 2034.699951171875, 452.75, 8445.75, -0.006563504014657771, 244.2920458514385)
    ,@("2024-09-18", 1879.449951171875, 1232.050048828125, 1889.400024414062, 2042.550048828125, 448.8500061035156, 8390.000091552734, -0.006600942302017657, 242.6794881519313)
    ,@("2024-09-19", 1866.650024414062, 1234.5, 1929.199951171875, 2025.699951171875, 450, 8406.049926757812, 0.00191297199403341, 243.1437272162923)
    ,@("2024-09-20", 1930.099975585938, 1246.550048828125, 1928.400024414062, 2062.300048828125, 446.7999877929688, 8507.750061035156, 0.01209844518691423, 246.0853882725606)
    ,@("2024-09-23", 2012.849975585938, 1229.900024414062, 1882.449951171875, 2039, 445.75, 8501.449951171875, -0.000740514215636784, 245.9031585442843)
    ,@("2024-09-24", 2026, 1238.800048828125, 1914.400024414062, 2064.550048828125, 439.25, 8561.500122070312, 0.007063521075032611, 247.640100687079)
    ,@("2024-09-25", 2051.39990234375, 1252.949951171875, 1910.849975585938, 2112.050048828125, 441.5499877929688, 8651.899841308594, 0.01055886444540762, 250.2548989414809)
    ,@("2024-09-26", 2134.800048828125, 1277.199951171875, 1857.550048828125, 2155.60009765625, 419.25, 8682.900146484375, 0.003583063343818423, 251.1515780964891)
    ,@("2024-09-27", 2119.5, 1275.75, 1866.699951171875, 2201.449951171875, 422.7999877929688, 8731.799865722656, 0.005631726544509473, 252.5659951055506)
)

$startRow = 651
for ($k = 0; $k -lt $newRows.Count; $k++) {
    $r = $startRow + $k
    $row = $newRows[$k]
    # Column A holds the date as literal text (matches the rest of the sheet,
    # which stores dates as plain "yyyy-MM-dd" strings, not Excel date serials).
    $aCell = $ws.Cells.Item($r, 1)
    $aCell.NumberFormat = "@"
    $aCell.Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $ws.Cells.Item($r, 8).Value = $row[6]
    $ws.Cells.Item($r, 9).Value = $row[7]
    $ws.Cells.Item($r, 10).Value = $row[8]
}
